{"js": "// Apply the nmap/meterpreter timestamp, latency & port updates described by\n// the diff. Uses body.search(...) + Range.insertText(..., Replace) so that\n// surrounding run/paragraph formatting is preserved.\n//\n// IMPORTANT: \"Host is up (0.00012s latency).\" becomes \"...0.00013s...\" while\n// the *other* host's \"Host is up (0.00013s latency).\" becomes \"...0.00017s...\".\n// If we searched/replaced sequentially, the second search for \"0.00013s\"\n// would also match the text we just wrote in the first replacement. To avoid\n// that, every search below is issued against the ORIGINAL document text\n// before any insertText() calls run, and all the resulting Range objects are\n// resolved/replaced together afterwards.\n\nconst body = context.document.body;\n\nconst searchNmapStart = body.search(\n  \"# Nmap 7.80 scan initiated Thu Jul 11 10:36:53 2024 as: nmap -sV -sC -Pn --script http-title -iL targets.txt -oN nmap_results.txt\",\n  { matchCase: true }\n);\nconst searchLatency225 = body.search(\"Host is up (0.00012s latency).\", { matchCase: true });\nconst searchLatency226 = body.search(\"Host is up (0.00013s latency).\", { matchCase: true });\nconst searchNmapDone = body.search(\n  \"# Nmap done at Thu Jul 11 10:37:01 2024 -- 2 IP addresses (2 hosts up) scanned in 7.99 seconds\",\n  { matchCase: true }\n);\nconst searchPort = body.search(\"10.33.102.225:40740\", { matchCase: true });\nconst searchMeterpreterTime = body.search(\"10:38:02\", { matchCase: true });\n\nsearchNmapStart.load(\"items\");\nsearchLatency225.load(\"items\");\nsearchLatency226.load(\"items\");\nsearchNmapDone.load(\"items\");\nsearchPort.load(\"items\");\nsearchMeterpreterTime.load(\"items\");\nawait context.sync();\n\nif (searchNmapStart.items.length !== 1) {\n  throw new Error(\"Expected 1 match for nmap start line, found \" + searchNmapStart.items.length);\n}\nif (searchLatency225.items.length !== 1) {\n  throw new Error(\"Expected 1 match for 10.33.102.225 latency line, found \" + searchLatency225.items.length);\n}\nif (searchLatency226.items.length !== 1) {\n  throw new Error(\"Expected 1 match for 10.33.102.226 latency line, found \" + searchLatency226.items.length);\n}\nif (searchNmapDone.items.length !== 1) {\n  throw new Error(\"Expected 1 match for nmap done line, found \" + searchNmapDone.items.length);\n}\nif (searchPort.items.length !== 2) {\n  throw new Error(\"Expected 2 matches for meterpreter port, found \" + searchPort.items.length);\n}\nif (searchMeterpreterTime.items.length !== 1) {\n  throw new Error(\"Expected 1 match for meterpreter session timestamp, found \" + searchMeterpreterTime.items.length);\n}\n\n// 1) First nmap scan line timestamp: 10:36:53 -> 10:40:18\nsearchNmapStart.items[0].insertText(\n  \"# Nmap 7.80 scan initiated Thu Jul 11 10:40:18 2024 as: nmap -sV -sC -Pn --script http-title -iL targets.txt -oN nmap_results.txt\",\n  Word.InsertLocation.replace\n);\n\n// 2) Latency for 10.33.102.225: 0.00012s -> 0.00013s\nsearchLatency225.items[0].insertText(\"Host is up (0.00013s latency).\", Word.InsertLocation.replace);\n\n// 3) Latency for 10.33.102.226: 0.00013s -> 0.00017s\nsearchLatency226.items[0].insertText(\"Host is up (0.00017s latency).\", Word.InsertLocation.replace);\n\n// 4) Nmap completion line: timestamp 10:37:01 -> 10:40:26, duration 7.99 -> 7.94 seconds\nsearchNmapDone.items[0].insertText(\n  \"# Nmap done at Thu Jul 11 10:40:26 2024 -- 2 IP addresses (2 hosts up) scanned in 7.94 seconds\",\n  Word.InsertLocation.replace\n);\n\n// 5) Meterpreter session port 40740 -> 33612 (appears twice: the \"session\n//    opened\" log line and the sessions table's wrapped connection column).\nfor (let i = 0; i < searchPort.items.length; i++) {\n  searchPort.items[i].insertText(\"10.33.102.225:33612\", Word.InsertLocation.replace);\n}\n\n// 6) Meterpreter session opened timestamp: 10:38:02 -> 10:41:21\nsearchMeterpreterTime.items[0].insertText(\"10:41:21\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Apply the nmap/meterpreter timestamp, latency & port updates described by\n# the diff, using Word COM Find/Replace (Range.Find.Execute).\n#\n# NOTE on the two \"Host is up (0.0001Xs latency).\" lines: the value for\n# 10.33.102.225 (0.00012s -> 0.00013s) and for 10.33.102.226\n# (0.00013s -> 0.00017s) collide textually after the first replace (both\n# would read \"0.00013s...\"). A fresh Range.Find.Execute always restarts\n# scanning from that range's Start, so re-searching the whole document for\n# \"0.00013s\" after the first edit would incorrectly re-match the text we\n# just wrote for .225 instead of the original .226 occurrence.\n#\n# To avoid that, both replacements reuse the SAME Range object ($nmapRng),\n# scoped to just the nmap-output paragraph. Word's Find naturally continues\n# from the end of the previous match when the same Range is reused, so the\n# second Find (for \"0.00013s\") resumes after the first replacement and only\n# finds the original .226 occurrence.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph holding the raw nmap console output (unique anchor\n# text), then scope a Range to it for the two ordered latency replacements.\n$nmapPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Text -like \"*Nmap 7.80 scan initiated*\") {\n        $nmapPara = $para\n        break\n    }\n}\nif ($nmapPara -eq $null) {\n    throw \"Could not locate the nmap output paragraph\"\n}\n$nmapRng = $nmapPara.Range\n\n# 1) First nmap scan line timestamp: 10:36:53 -> 10:40:18\n$ok = $nmapRng.Find.Execute(\n    \"# Nmap 7.80 scan initiated Thu Jul 11 10:36:53 2024 as: nmap -sV -sC -Pn --script http-title -iL targets.txt -oN nmap_results.txt\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"# Nmap 7.80 scan initiated Thu Jul 11 10:40:18 2024 as: nmap -sV -sC -Pn --script http-title -iL targets.txt -oN nmap_results.txt\",\n    1\n)\nif (-not $ok) { throw \"Failed to update nmap start timestamp\" }\n\n# 2) Latency for 10.33.102.225: 0.00012s -> 0.00013s (first occurrence from\n#    the current, advanced search position)\n$ok = $nmapRng.Find.Execute(\n    \"Host is up (0.00012s latency).\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Host is up (0.00013s latency).\",\n    1\n)\nif (-not $ok) { throw \"Failed to update 10.33.102.225 latency\" }\n\n# 3) Latency for 10.33.102.226: 0.00013s -> 0.00017s (continues searching\n#    forward from where step 2 left off, so only the original .226 text --\n#    not the one we just wrote -- is matched)\n$ok = $nmapRng.Find.Execute(\n    \"Host is up (0.00013s latency).\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Host is up (0.00017s latency).\",\n    1\n)\nif (-not $ok) { throw \"Failed to update 10.33.102.226 latency\" }\n\n# 4) Nmap completion line: timestamp 10:37:01 -> 10:40:26, duration 7.99 -> 7.94 seconds\n$ok = $nmapRng.Find.Execute(\n    \"# Nmap done at Thu Jul 11 10:37:01 2024 -- 2 IP addresses (2 hosts up) scanned in 7.99 seconds\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"# Nmap done at Thu Jul 11 10:40:26 2024 -- 2 IP addresses (2 hosts up) scanned in 7.94 seconds\",\n    1\n)\nif (-not $ok) { throw \"Failed to update nmap done line\" }\n\n# 5) Meterpreter session port 40740 -> 33612 (appears twice: the \"session\n#    opened\" log line and the sessions table's wrapped connection column).\n#    Both occurrences get the same replacement text, so a document-wide\n#    Replace All is unambiguous.\n$allRng = $d.Content\n$ok = $allRng.Find.Execute(\n    \"10.33.102.225:40740\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"10.33.102.225:33612\",\n    2\n)\nif (-not $ok) { throw \"Failed to update meterpreter session port\" }\n\n# 6) Meterpreter session opened timestamp: 10:38:02 -> 10:41:21\n$allRng = $d.Content\n$ok = $allRng.Find.Execute(\n    \"10:38:02\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"10:41:21\",\n    1\n)\nif (-not $ok) { throw \"Failed to update meterpreter session timestamp\" }\n"}
